$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("scenarios")

# Insert a new row above the current row 7 (StartShutCost), pushing the
# existing rows 7-9 down to 8-10, and fill it with the new
# "RollingFixDispatch" scenario entries (mirrors the A/H/I pattern already
# used for the "RollingHorizon" row).
$ws.Rows("7:7").Insert()
$ws.Range("A7").Value = "RollingFixDispatch"
$ws.Range("H7").Value = "RollingFixDispatch"
$ws.Range("I7").Value = "RollingFixDispatch"

# Column A needs to widen a bit to fit the new, longer label.
$ws.Columns("A").ColumnWidth = 16.67

# Make "scenarios" the active/selected sheet (it was "model_config" before).
$ws.Activate()
$ws.Range("C13").Select()
